$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain text (their new values are plain
# numeric-looking strings that Excel would otherwise auto-convert to numbers).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '59.747.18'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '2.402.59'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '550.41'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  +3.52%  '
$ws.Range('E9').Value = '  -2.21%  '
$ws.Range('D10').Value = '5.69'
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  -2.07%  '
$ws.Range('D13').Value = '25.30'
$ws.Range('E13').Value = '  +2.43%  '
$ws.Range('D14').Value = '2.830.58'
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').Value = '59.675.98'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').Value = '2.433.30'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '11.32'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('E19').Value = '  -1.28%  '
$ws.Range('D20').Value = '328.70'
$ws.Range('E20').Value = '  -1.95%  '
$ws.Range('E21').Value = '  -4.32%  '
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('E23').Value = '  +2.93%  '
$ws.Range('D24').Value = '0.173'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').Value = '1.36'
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('E29').Value = '  -2.28%  '
$ws.Range('D30').Value = '168.42'
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('D31').Value = '6.05'
$ws.Range('E31').Value = '  -3.90%  '
$ws.Range('D32').Value = '18.61'
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('E33').Value = '  -1.86%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -2.31%  '
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').Value = '319.73'
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('D40').Value = '0.406'
$ws.Range('E40').Value = '  -3.39%  '
$ws.Range('E41').Value = '  -2.39%  '
$ws.Range('D42').Value = '139.08'
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range('D43').Value = '0.0968'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = '19.58'
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('D45').Value = '0.0513'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').Value = '0.577'
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('D48').Value = '0.385'
$ws.Range('E48').Value = '  -6.53%  '
$ws.Range('D49').Value = '17.53'
$ws.Range('E49').Value = '  -2.52%  '
$ws.Range('D50').Value = '11.05'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  -3.33%  '
